$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.717.13'
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").Value = '3.256.53'
$ws.Range("E3").Value = '  +2.35%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '607.42'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.85%  '

$ws.Range("E6").Value = '  +2.04%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.255.01'
$ws.Range("E8").Value = '  +2.31%  '

$ws.Range("E9").Value = '  -0.37%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.162'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.86%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.94'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +6.03%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.508'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.98%  '

$ws.Range("E13").Value = '  +1.47%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '39.30'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.54%  '

$ws.Range("D15").Value = '3.787.52'
$ws.Range("E15").Value = '  +2.20%  '

$ws.Range("D16").Value = '66.732.01'
$ws.Range("E16").Value = '  +0.31%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '7.44'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.21%  '

$ws.Range("D18").Value = '3.252.94'
$ws.Range("E18").Value = '  +2.21%  '

$ws.Range("E19").Value = '  +1.16%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '506.92'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.26%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '15.44'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.17%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.753'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +2.80%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '8.14'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.08%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '14.85'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.30%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '86.66'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.25%  '

$ws.Range("E26").Value = '  +0.10%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.145'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +62.70%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '3.03'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.13%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.08'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.24%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.40'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.10%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '2.89'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -7.03%  '

$ws.Range("E32").Value = '  -2.11%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '28.14'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.15%  '

$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("E35").Value = '  -4.13%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '6.45'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.35%  '

$ws.Range("D37").Value = '0.0₃0796'
$ws.Range("E37").Value = '  +16.82%  '

$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.35'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +19.90%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '55.64'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.43%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '494.60'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -4.00%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.0428'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.85%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.128'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("E44").Value = '  -1.71%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.49'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +2.30%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.969.35'
$ws.Range("E46").Value = '  +4.36%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '28.75'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.21%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.49'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +4.07%  '

$ws.Range("E49").Value = '  +2.71%  '

$ws.Range("E50").Value = '  -0.05%  '

$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '121.49'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.70%  '
